# Add average attack distances to the "attacks" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect its contents.
$ws.Name = "all_binned_counts"

# New column header.
$ws.Range("J1").Value = "avg_distance_km"

# Refresh header formatting across the whole header row (bold, centered,
# top-aligned, thin box border) so the new column matches the rest.
$headerRange = $ws.Range("A1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

# Average attack distance (km) per period, newly computed column.
$avgDistance = @{
    2  = 3.76326812791891
    3  = 4.385888662670677
    4  = 4.11879765724225
    5  = 4.897854013092467
    6  = 4.476229029663664
    7  = 3.876734017006266
    8  = 3.709505535754547
    9  = 4.299951632670066
    10 = 4.675176933039502
    11 = 4.286881911154486
    12 = 4.541151156676015
    13 = 4.723351366358533
    14 = 4.561800693261488
    15 = 4.508389709238768
    16 = 4.518039712920642
    17 = 4.490206231553414
    18 = 4.171829734708102
    19 = 4.136376498618044
}

foreach ($row in 2..19) {
    $ws.Range("J$row").Value = $avgDistance[$row]
}

# Page margins reset to Excel's standard defaults (inches: 0.75/0.75/1/1/0.5/0.5).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
